$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("N47").Value = $null

$ws.Range("H92").Value = 361
$ws.Range("I92").Value = 381.6
$ws.Range("J92").Value = 155
$ws.Range("K92").Value = 381.6
$ws.Range("L92").Value = 155
$ws.Range("M92").Value = 866.4
$ws.Range("N92").Value = -2651

$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492

$ws.Range("H105").Value = 11200
$ws.Range("I105").Value = 11200
$ws.Range("K105").Value = 11200
$ws.Range("M105").Value = -7706

$ws.Range("H116").Value = 6208.25
$ws.Range("J116").Value = 6062.5
$ws.Range("L116").Value = 6062.5
$ws.Range("N116").Value = -12946.5

$ws.Range("H131").Value = 410
$ws.Range("I131").Value = 410
$ws.Range("K131").Value = 1230
$ws.Range("M131").Value = 3810

$ws.Range("H132").Value = 3831.4211
$ws.Range("J132").Value = 4911.1113
$ws.Range("L132").Value = 14733.3339
$ws.Range("N132").Value = -19793.3339

$ws.Range("H138").Value = 9953.541999999999
$ws.Range("I138").Value = 7587.4
$ws.Range("J138").Value = 10576.211
$ws.Range("K138").Value = 22762.2
$ws.Range("L138").Value = 31728.633
$ws.Range("M138").Value = -17622.2
$ws.Range("N138").Value = -42008.633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8176.914
$ws.Range("I32").Value = 6686.6562
$ws.Range("K32").Value = 6686.6562
$ws.Range("M32").Value = -6399.6562

$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

$ws.Range("H132").Value = 1418
$ws.Range("I132").Value = 1418
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4254
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1724
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 11875
$ws.Range("J100").Value = 11875
$ws.Range("L100").Value = 11875
$ws.Range("N100").Value = -14039

$ws.Range("H103").Value = 7999
$ws.Range("J103").Value = 7999
$ws.Range("L103").Value = 7999
$ws.Range("N103").Value = -10343

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2811.1333
$ws.Range("I31").Value = 2819.6155
$ws.Range("J31").Value = 2756
$ws.Range("K31").Value = 2819.6155
$ws.Range("L31").Value = 2756
$ws.Range("M31").Value = -2524.6155
$ws.Range("N31").Value = -3346

$ws.Range("H34").Value = 2811.1333
$ws.Range("I34").Value = 2819.6155
$ws.Range("J34").Value = 2756
$ws.Range("K34").Value = 2819.6155
$ws.Range("L34").Value = 2756
$ws.Range("M34").Value = -2617.6155
$ws.Range("N34").Value = -3160

$ws.Range("H99").Value = 9305.416999999999
$ws.Range("I99").Value = 8612.666999999999
$ws.Range("K99").Value = 8612.666999999999
$ws.Range("M99").Value = -7114.666999999999

$ws.Range("H106").Value = 200223.33
$ws.Range("J106").Value = 200223.33
$ws.Range("L106").Value = 200223.33
$ws.Range("N106").Value = -202747.33

$ws.Range("H122").Value = 1519.8
$ws.Range("I122").Value = 1149.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3449.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -999.25
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 9305.416999999999
$ws.Range("I126").Value = 8612.666999999999
$ws.Range("K126").Value = 25838.001
$ws.Range("M126").Value = -23368.001

$ws.Range("H132").Value = 6802.154
$ws.Range("I132").Value = 2491.0625
$ws.Range("K132").Value = 7473.1875
$ws.Range("M132").Value = -4943.1875

$ws.Range("H134").Value = 3556.077
$ws.Range("I134").Value = 3338.5908
$ws.Range("J134").Value = 4752.25
$ws.Range("K134").Value = 10015.7724
$ws.Range("L134").Value = 14256.75
$ws.Range("M134").Value = -7480.7724
$ws.Range("N134").Value = -19326.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1304.4445
$ws.Range("I68").Value = 1304.4445
$ws.Range("K68").Value = 3913.3335
$ws.Range("M68").Value = -3102.3335

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null

$ws.Range("H71").Value = 1304.4445
$ws.Range("I71").Value = 1304.4445
$ws.Range("K71").Value = 11740.0005
$ws.Range("M71").Value = -7684.0005

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H122").Value = 1354.9412
$ws.Range("I122").Value = 1616.25
$ws.Range("J122").Value = 1274.5385
$ws.Range("K122").Value = 14546.25
$ws.Range("L122").Value = 11470.8465
$ws.Range("M122").Value = -12096.25
$ws.Range("N122").Value = -16370.8465

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1200
$ws.Range("I102").Value = 1200
$ws.Range("K102").Value = 1200
$ws.Range("M102").Value = 422

$ws.Range("H132").Value = 5970.6665
$ws.Range("I132").Value = 4456
$ws.Range("K132").Value = 13368
$ws.Range("M132").Value = -10838

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2744
$ws.Range("J61").Value = 3005
$ws.Range("L61").Value = 3005
$ws.Range("N61").Value = -3409

$ws.Range("H113").Value = 2744
$ws.Range("J113").Value = 3005
$ws.Range("L113").Value = 3005
$ws.Range("N113").Value = -7345

$ws.Range("H122").Value = 8993.333000000001
$ws.Range("J122").Value = 8988
$ws.Range("L122").Value = 26964
$ws.Range("N122").Value = -31864

$ws.Range("H136").Value = 4966.6665
$ws.Range("I136").Value = 4450
$ws.Range("K136").Value = 13350
$ws.Range("M136").Value = -10800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 804.25
$ws.Range("J2").Value = 1597.5
$ws.Range("L2").Value = 1597.5
$ws.Range("N2").Value = -1821.5

$ws.Range("H4").Value = 4825.75
$ws.Range("J4").Value = 4101
$ws.Range("L4").Value = 4101
$ws.Range("N4").Value = -4327

$ws.Range("H75").Value = 22997.8
$ws.Range("J75").Value = 22997.8
$ws.Range("L75").Value = 22997.8
$ws.Range("N75").Value = -24869.8

$ws.Range("H78").Value = 22997.8
$ws.Range("J78").Value = 22997.8
$ws.Range("L78").Value = 68993.39999999999
$ws.Range("N78").Value = -78353.39999999999
